$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (row index -> values for columns A..T)
# A/B/C/D are text labels; E..T are numeric values
$rows = @(
    @{ RowNum=2;  A="ECs";  B="Cd44"; C="Sele"; D="ECs";  E=3; F=1; G=31.82741333333333;  H=95.48223999999999;  I=0.114390792932228;  J=0.114390792932228;  K=3; L=1;                   M=4.699506666666667;  N=14.09852;  O=0.9660495246229048; P=0.9660495246229047; Q=149.5731411427556;  R=1346.1582702848;    S=0.110507171133416;    T=0.110507171133416 },
    @{ RowNum=3;  A="ECs";  B="Cd44"; C="Sele"; D="FAPs"; E=3; F=1; G=31.82741333333333;  H=95.48223999999999;  I=0.114390792932228;  J=0.114390792932228;  K=1; L=0.3333333333333333;  M=0.1651576666666667; N=0.495473;  O=0.03395047537709522; P=0.03395047537709522; Q=5.256541322168888;  R=47.30887189951999;  S=0.003883621798812005; T=0.003883621798812006 },
    @{ RowNum=4;  A="FAPs"; B="Cd44"; C="Sele"; D="ECs";  E=3; F=1; G=85.46317833333335;  H=256.389535;         I=0.307162904935779;  J=0.307162904935779;  K=3; L=1;                   M=4.699506666666667;  N=14.09852;  O=0.9660495246229048; P=0.9660495246229047; Q=401.6347763320223;  R=3614.712986988201;  S=0.2967345782949998;   T=0.2967345782949998 },
    @{ RowNum=5;  A="FAPs"; B="Cd44"; C="Sele"; D="FAPs"; E=3; F=1; G=85.46317833333335;  H=256.389535;         I=0.307162904935779;  J=0.307162904935779;  K=1; L=0.3333333333333333;  M=0.1651576666666667; N=0.495473;  O=0.03395047537709522; P=0.03395047537709522; Q=14.11489911945056;  R=127.034092075055;   S=0.0104283266407792;   T=0.0104283266407792 },
    @{ RowNum=6;  A="M2";   B="Cd44"; C="Sele"; D="ECs";  E=3; F=1; G=122.2478306666667; H=366.743492;         I=0.4393704929064738; J=0.4393704929064738; K=3; L=1;                   M=4.699506666666667;  N=14.09852;  O=0.9660495246229048; P=0.9660495246229047; Q=574.5044952035379;  R=5170.54045683184;   S=0.4244536558056304;   T=0.4244536558056304 },
    @{ RowNum=7;  A="M2";   B="Cd44"; C="Sele"; D="FAPs"; E=3; F=1; G=122.2478306666667; H=366.743492;         I=0.4393704929064738; J=0.4393704929064738; K=1; L=0.3333333333333333;  M=0.1651576666666667; N=0.495473;  O=0.03395047537709522; P=0.03395047537709522; Q=20.19016646796845;  R=181.711498211716;   S=0.01491683710084343;  T=0.01491683710084343 },
    @{ RowNum=8;  A="sCs";  B="Cd44"; C="Sele"; D="ECs";  E=3; F=1; G=38.69562533333333; H=116.086876;         I=0.1390758092255191; J=0.1390758092255191; K=3; L=1;                   M=4.699506666666667;  N=14.09852;  O=0.9660495246229048; P=0.9660495246229047; Q=181.8503492248356;  R=1636.65314302352;   S=0.1343541193888586;   T=0.1343541193888586 },
    @{ RowNum=9;  A="sCs";  B="Cd44"; C="Sele"; D="FAPs"; E=3; F=1; G=38.69562533333333; H=116.086876;         I=0.1390758092255191; J=0.1390758092255191; K=1; L=0.3333333333333333;  M=0.1651576666666667; N=0.495473;  O=0.03395047537709522; P=0.03395047537709522; Q=6.390879190260889;  R=57.51791271234799;  S=0.004721689836660579; T=0.004721689836660579 }
)

foreach ($row in $rows) {
    $r = $row.RowNum
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}
